$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Alinea "1g" (row 14): progress updated from 0 to 80, and an owner ("Eduardo") assigned.
# The "Status" column (F) recalculates automatically via its existing formula.
$ws.Range("D14").Value = 80
$ws.Range("E14").Value = "Eduardo"

# Alinea "1n" (row 21): owner assigned ("Eduardo"); progress remains 0 (still TODO).
$ws.Range("E21").Value = "Eduardo"

# Move/select the active cell to K15, as in the saved workbook.
$ws.Range("K15").Select()
